$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6114.3335
$ws.Range("I64").Value = 3519.2307
$ws.Range("J64").Value = 9181.272000000001
$ws.Range("K64").Value = 3519.2307
$ws.Range("L64").Value = 9181.272000000001
$ws.Range("M64").Value = -3271.2307
$ws.Range("N64").Value = -9677.272000000001
$ws.Range("H67").Value = 6114.3335
$ws.Range("I67").Value = 3519.2307
$ws.Range("J67").Value = 9181.272000000001
$ws.Range("K67").Value = 3519.2307
$ws.Range("L67").Value = 9181.272000000001
$ws.Range("M67").Value = -2661.2307
$ws.Range("N67").Value = -10897.272
$ws.Range("H137").Value = 2038.325
$ws.Range("I137").Value = 1916.8064
$ws.Range("J137").Value = 2456.889
$ws.Range("K137").Value = 5750.4192
$ws.Range("L137").Value = 7370.667
$ws.Range("M137").Value = -3200.4192
$ws.Range("N137").Value = -12470.667
$ws.Range("H138").Value = 6413691.5
$ws.Range("I138").Value = 2254.1667
$ws.Range("J138").Value = 7579407.5
$ws.Range("K138").Value = 6762.500100000001
$ws.Range("L138").Value = 22738222.5
$ws.Range("M138").Value = -1622.500100000001
$ws.Range("N138").Value = -22748502.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1645.6154
$ws.Range("J2").Value = 1271
$ws.Range("L2").Value = 1271
$ws.Range("N2").Value = -1497
$ws.Range("H45").Value = 6938.3335
$ws.Range("I45").Value = 7826.067
$ws.Range("K45").Value = 7826.067
$ws.Range("M45").Value = -7449.067
$ws.Range("H53").Value = 17990
$ws.Range("I53").Value = 17990
$ws.Range("K53").Value = 17990
$ws.Range("M53").Value = -17308
$ws.Range("H61").Value = 5253.1
$ws.Range("I61").Value = 4504.4287
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 4504.4287
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -4292.4287
$ws.Range("N61").Value = -7424
$ws.Range("H74").Value = 11531.95
$ws.Range("I74").Value = 2730
$ws.Range("J74").Value = 32069.834
$ws.Range("K74").Value = 2730
$ws.Range("L74").Value = 32069.834
$ws.Range("M74").Value = -1856
$ws.Range("N74").Value = -33817.834
$ws.Range("H77").Value = 11531.95
$ws.Range("I77").Value = 2730
$ws.Range("J77").Value = 32069.834
$ws.Range("K77").Value = 13650
$ws.Range("L77").Value = 160349.17
$ws.Range("M77").Value = -9282
$ws.Range("N77").Value = -169085.17
$ws.Range("H102").Value = 3026.4138
$ws.Range("I102").Value = 3185.3462
$ws.Range("K102").Value = 3185.3462
$ws.Range("M102").Value = -1563.3462
$ws.Range("H116").Value = 1645.6154
$ws.Range("J116").Value = 1271
$ws.Range("L116").Value = 1271
$ws.Range("N116").Value = -5859
$ws.Range("H122").Value = 1999.5834
$ws.Range("I122").Value = 1829.6
$ws.Range("K122").Value = 5488.799999999999
$ws.Range("M122").Value = -3038.799999999999
$ws.Range("H128").Value = 80000
$ws.Range("J128").Value = 80000
$ws.Range("L128").Value = 80000
$ws.Range("N128").Value = -89960
$ws.Range("H136").Value = 5253.1
$ws.Range("I136").Value = 4504.4287
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 13513.2861
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -10963.2861
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1645.6154
$ws.Range("J3").Value = 1271
$ws.Range("L3").Value = 1271
$ws.Range("N3").Value = -1499
$ws.Range("H86").Value = 4577
$ws.Range("I86").Value = 5521.5713
$ws.Range("K86").Value = 5521.5713
$ws.Range("M86").Value = -4398.5713
$ws.Range("H89").Value = 4577
$ws.Range("I89").Value = 5521.5713
$ws.Range("K89").Value = 27607.8565
$ws.Range("M89").Value = -21991.8565
$ws.Range("H99").Value = 50444.117
$ws.Range("I99").Value = 45289.957
$ws.Range("K99").Value = 45289.957
$ws.Range("M99").Value = -43791.957
$ws.Range("H134").Value = 2801.2222
$ws.Range("I134").Value = 2666
$ws.Range("K134").Value = 7998
$ws.Range("M134").Value = -5463

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1794.92
$ws.Range("I31").Value = 1591.4783
$ws.Range("J31").Value = 4134.5
$ws.Range("K31").Value = 1591.4783
$ws.Range("L31").Value = 4134.5
$ws.Range("M31").Value = -1296.4783
$ws.Range("N31").Value = -4724.5
$ws.Range("H34").Value = 1794.92
$ws.Range("I34").Value = 1591.4783
$ws.Range("J34").Value = 4134.5
$ws.Range("K34").Value = 1591.4783
$ws.Range("L34").Value = 4134.5
$ws.Range("M34").Value = -1389.4783
$ws.Range("N34").Value = -4538.5
$ws.Range("H55").Value = 32371.2
$ws.Range("I55").Value = 7776
$ws.Range("J55").Value = 38520
$ws.Range("K55").Value = 7776
$ws.Range("L55").Value = 38520
$ws.Range("M55").Value = -7461
$ws.Range("N55").Value = -39150
$ws.Range("H88").Value = 17969.143
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 17969.143
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 17969.143
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -18781.143
$ws.Range("H91").Value = 17969.143
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 17969.143
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 17969.143
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -20777.143
$ws.Range("H99").Value = 3111.111
$ws.Range("I99").Value = 2928.5715
$ws.Range("J99").Value = 3750
$ws.Range("K99").Value = 2928.5715
$ws.Range("L99").Value = 3750
$ws.Range("M99").Value = -1430.5715
$ws.Range("N99").Value = -6746
$ws.Range("H122").Value = 3121.3333
$ws.Range("I122").Value = 2548.6667
$ws.Range("J122").Value = 4266.6665
$ws.Range("K122").Value = 7646.000100000001
$ws.Range("L122").Value = 12799.9995
$ws.Range("M122").Value = -5196.000100000001
$ws.Range("N122").Value = -17699.9995
$ws.Range("H126").Value = 3111.111
$ws.Range("I126").Value = 2928.5715
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 8785.7145
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -6315.7145
$ws.Range("N126").Value = -16190

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 369.82144
$ws.Range("J12").Value = 407.5263
$ws.Range("L12").Value = 1222.5789
$ws.Range("N12").Value = -1568.5789
$ws.Range("H132").Value = 1981.5714
$ws.Range("I132").Value = 1437.6316
$ws.Range("J132").Value = 2627.5
$ws.Range("K132").Value = 12938.6844
$ws.Range("L132").Value = 23647.5
$ws.Range("M132").Value = -10408.6844
$ws.Range("N132").Value = -28707.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 23000
$ws.Range("I55").Value = 12500
$ws.Range("K55").Value = 12500
$ws.Range("M55").Value = -12173
$ws.Range("H113").Value = 1465
$ws.Range("I113").Value = 1518
$ws.Range("K113").Value = 1518
$ws.Range("M113").Value = 652
$ws.Range("H122").Value = 4128.5713
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 4650
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 13950
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -18850

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 71998.664
$ws.Range("I23").Value = 56398.4
$ws.Range("K23").Value = 56398.4
$ws.Range("M23").Value = -56168.4
$ws.Range("H48").Value = 6999.6665
$ws.Range("I48").Value = 6999.6665
$ws.Range("K48").Value = 6999.6665
$ws.Range("M48").Value = -6338.6665
$ws.Range("H136").Value = 3559.1875
$ws.Range("I136").Value = 3253.2415
$ws.Range("J136").Value = 6516.6665
$ws.Range("K136").Value = 9759.7245
$ws.Range("L136").Value = 19549.9995
$ws.Range("M136").Value = -7209.7245
$ws.Range("N136").Value = -24649.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 952
$ws.Range("I100").Value = 990
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 1980
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -1439
$ws.Range("N100").Value = -2682
$ws.Range("H122").Value = 2654.1292
$ws.Range("I122").Value = 2415.52
$ws.Range("J122").Value = 3648.3333
$ws.Range("K122").Value = 7246.559999999999
$ws.Range("L122").Value = 10944.9999
$ws.Range("M122").Value = -4796.559999999999
$ws.Range("N122").Value = -15844.9999
